$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FTEO")

$ws.Range("D8").Value = 110400
$ws.Range("E8").Value = 101300
$ws.Range("F8").Value = 95400
$ws.Range("G8").Value = 56700
$ws.Range("H8").Value = 37700
$ws.Range("I8").Value = 42300
$ws.Range("J8").Value = 46400
$ws.Range("D9").Value = 63200
$ws.Range("E9").Value = 62100
$ws.Range("F9").Value = 52900
$ws.Range("G9").Value = 28400
$ws.Range("H9").Value = 20900
$ws.Range("I9").Value = 16500
$ws.Range("J9").Value = 14700
$ws.Range("D10").Value = 47300
$ws.Range("E10").Value = 39200
$ws.Range("F10").Value = 42500
$ws.Range("G10").Value = 28300
$ws.Range("H10").Value = 16800
$ws.Range("I10").Value = 25800
$ws.Range("J10").Value = 31800
$ws.Range("E12").Value = 800
$ws.Range("D14").Value = 7600
$ws.Range("D15").Value = 4500
$ws.Range("E15").Value = 4100
$ws.Range("D17").Value = 116500
$ws.Range("E17").Value = 113000
$ws.Range("F17").Value = 95100
$ws.Range("G17").Value = 54300
$ws.Range("H17").Value = 43400
$ws.Range("I17").Value = 34000
$ws.Range("J17").Value = 25000
$ws.Range("D18").Value = -6000
$ws.Range("E18").Value = -11700
$ws.Range("H18").Value = -5700
$ws.Range("I18").Value = 8300
$ws.Range("J18").Value = 21400
$ws.Range("D20").Value = -1300
$ws.Range("J20").Value = 100
$ws.Range("D21").Value = 3500
$ws.Range("F21").Value = 9100
$ws.Range("G21").Value = 9400
$ws.Range("H21").Value = -2000
$ws.Range("I21").Value = 10600
$ws.Range("J21").Value = 22900
$ws.Range("D23").Value = -7800
$ws.Range("E23").Value = -12400
$ws.Range("H23").Value = -6000
$ws.Range("I23").Value = 7900
$ws.Range("J23").Value = 21500
$ws.Range("E24").Value = -3900
$ws.Range("G24").Value = 1600
$ws.Range("J24").Value = 8800
$ws.Range("D26").Value = -7400
$ws.Range("E26").Value = -8500
$ws.Range("G26").Value = 2400
$ws.Range("H26").Value = -5400
$ws.Range("I26").Value = 4600
$ws.Range("J26").Value = 21500
$ws.Range("D27").Value = -7500
$ws.Range("E27").Value = -8600
$ws.Range("F27").Value = -1800
$ws.Range("G27").Value = 2400
$ws.Range("H27").Value = -5500
$ws.Range("I27").Value = 4500
$ws.Range("J27").Value = 12400
$ws.Range("D32").Value = 1300
$ws.Range("J32").Value = -100
$ws.Range("D33").Value = -7500
$ws.Range("E33").Value = -8600
$ws.Range("F33").Value = -1800
$ws.Range("G33").Value = 2400
$ws.Range("H33").Value = -5500
$ws.Range("I33").Value = 4500
$ws.Range("J33").Value = 12400
$ws.Range("D35").Value = -7500
$ws.Range("E35").Value = -8600
$ws.Range("F35").Value = -1800
$ws.Range("G35").Value = 2400
$ws.Range("H35").Value = -5500
$ws.Range("I35").Value = 4500
$ws.Range("J35").Value = 12400
$ws.Range("D41").Value = 46400
$ws.Range("E41").Value = 41000
$ws.Range("F41").Value = 16300
$ws.Range("G41").Value = 24600
$ws.Range("H41").Value = 13400
$ws.Range("I41").Value = 11800
$ws.Range("J41").Value = 21800
$ws.Range("D43").Value = 23100
$ws.Range("E43").Value = 28500
$ws.Range("F43").Value = 25300
$ws.Range("G43").Value = 14000
$ws.Range("H43").Value = 7900
$ws.Range("I43").Value = 10600
$ws.Range("J43").Value = 9100
$ws.Range("E45").Value = 11500
$ws.Range("F45").Value = 7300
$ws.Range("H45").Value = 3600
$ws.Range("I45").Value = 2700
$ws.Range("J45").Value = 2300
$ws.Range("D46").Value = 73000
$ws.Range("E46").Value = 81700
$ws.Range("F46").Value = 48900
$ws.Range("G46").Value = 41800
$ws.Range("H46").Value = 24900
$ws.Range("I46").Value = 25000
$ws.Range("J46").Value = 33200
$ws.Range("D47").Value = 7700
$ws.Range("E47").Value = 4500
$ws.Range("F47").Value = 5800
$ws.Range("G47").Value = 4600
$ws.Range("H47").Value = 3400
$ws.Range("J47").Value = 2500
$ws.Range("D48").Value = 13800
$ws.Range("E48").Value = 10600
$ws.Range("F48").Value = 10000
$ws.Range("G48").Value = 7100
$ws.Range("H48").Value = 6100
$ws.Range("I48").Value = 6800
$ws.Range("D49").Value = 39600
$ws.Range("E49").Value = 45400
$ws.Range("F49").Value = 48500
$ws.Range("G49").Value = 14200
$ws.Range("H49").Value = 8500
$ws.Range("I49").Value = 7600
$ws.Range("J49").Value = 4400
$ws.Range("E52").Value = 3800
$ws.Range("D54").Value = 130800
$ws.Range("E54").Value = 146100
$ws.Range("F54").Value = 116800
$ws.Range("G54").Value = 69100
$ws.Range("H54").Value = 44200
$ws.Range("I54").Value = 43000
$ws.Range("J54").Value = 44100
$ws.Range("F57").Value = 3300
$ws.Range("H57").Value = 600
$ws.Range("D58").Value = 22200
$ws.Range("E58").Value = 16300
$ws.Range("F58").Value = 7000
$ws.Range("G58").Value = 4500
$ws.Range("H58").Value = 2500
$ws.Range("I58").Value = 1700
$ws.Range("J58").Value = 1900
$ws.Range("D59").Value = 13900
$ws.Range("E59").Value = 12500
$ws.Range("F59").Value = 21400
$ws.Range("G59").Value = 7800
$ws.Range("H59").Value = 4100
$ws.Range("I59").Value = 6900
$ws.Range("J59").Value = 13400
$ws.Range("D60").Value = 37400
$ws.Range("E60").Value = 33400
$ws.Range("F60").Value = 31700
$ws.Range("G60").Value = 13100
$ws.Range("H60").Value = 7100
$ws.Range("I60").Value = 9300
$ws.Range("J60").Value = 15800
$ws.Range("D61").Value = 47900
$ws.Range("E61").Value = 55000
$ws.Range("F61").Value = 34500
$ws.Range("G61").Value = 7000
$ws.Range("H61").Value = 4000
$ws.Range("I61").Value = 4000
$ws.Range("D62").Value = 6100
$ws.Range("E62").Value = 12200
$ws.Range("F62").Value = 8500
$ws.Range("H62").Value = 1200
$ws.Range("D66").Value = 91800
$ws.Range("E66").Value = 101000
$ws.Range("F66").Value = 74900
$ws.Range("G66").Value = 22000
$ws.Range("H66").Value = 12400
$ws.Range("I66").Value = 14300
$ws.Range("J66").Value = 20300
$ws.Range("D72").Value = -4200
$ws.Range("E72").Value = 2800
$ws.Range("F72").Value = 12900
$ws.Range("G72").Value = 14800
$ws.Range("H72").Value = 12300
$ws.Range("I72").Value = 18300
$ws.Range("J72").Value = 14600
$ws.Range("D76").Value = 39000
$ws.Range("E76").Value = 45100
$ws.Range("F76").Value = 41900
$ws.Range("G76").Value = 47100
$ws.Range("H76").Value = 31700
$ws.Range("I76").Value = 28700
$ws.Range("J76").Value = 23900
$ws.Range("D81").Value = -7500
$ws.Range("E81").Value = -8600
$ws.Range("F81").Value = -1800
$ws.Range("G81").Value = 2400
$ws.Range("H81").Value = -5500
$ws.Range("I81").Value = 4500
$ws.Range("J81").Value = 12400
$ws.Range("D83").Value = 10800
$ws.Range("E83").Value = 10600
$ws.Range("F83").Value = 9000
$ws.Range("G83").Value = 5400
$ws.Range("H83").Value = 4000
$ws.Range("D89").Value = 13500
$ws.Range("F89").Value = 5600
$ws.Range("G89").Value = 9200
$ws.Range("J89").Value = 18800
$ws.Range("D91").Value = -1900
$ws.Range("E91").Value = -3000
$ws.Range("F91").Value = -6100
$ws.Range("G91").Value = -1800
$ws.Range("H91").Value = -1600
$ws.Range("I91").Value = -4600
$ws.Range("J91").Value = -2200
$ws.Range("D94").Value = -6400
$ws.Range("E94").Value = -14900
$ws.Range("F94").Value = -42600
$ws.Range("G94").Value = -12800
$ws.Range("H94").Value = -5900
$ws.Range("I94").Value = -10300
$ws.Range("J94").Value = -5000
$ws.Range("F96").Value = -1000
$ws.Range("D100").Value = -1300
$ws.Range("E100").Value = 41300
$ws.Range("F100").Value = 29100
$ws.Range("G100").Value = 15200
$ws.Range("H100").Value = 6800
$ws.Range("I100").Value = -2600
$ws.Range("J100").Value = 2000
$ws.Range("D102").Value = 5400
$ws.Range("E102").Value = 24700
$ws.Range("F102").Value = -8300
$ws.Range("G102").Value = 12100
$ws.Range("H102").Value = 1700
$ws.Range("I102").Value = -11000
$ws.Range("J102").Value = 15700
